# Updates the cryptos list with freshly scraped price/volume data.
# Most rows keep the same coin in place and only refresh Price (D) and
# Volume(1h) (E). A few rows swap which coin occupies which row (the
# B/C/D/E cells for those row-pairs are fully replaced).
#
# The Price column stores plain decimal-looking text (e.g. "263.99").
# Assigning such a string straight to .Value would make Excel infer a
# genuine number and silently reformat it (263.99000000000001 due to
# float rounding) instead of keeping the literal text from the feed. To
# preserve the exact original text we prefix values that parse as plain
# numbers with a leading apostrophe, which is how Excel is told "store
# this as text" without touching the cell's NumberFormat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# --- Rows whose coin stays the same; only Price / Volume(1h) change ---
$priceVolumeUpdates = @(
    @{ Row = 2;  Price = "43.156.59";  Volume = "  -1.52%  " },
    @{ Row = 3;  Price = "2.273.72";   Volume = "  -1.02%  " },
    @{ Row = 4;  Price = $null;        Volume = "  -0.22%  " },
    @{ Row = 5;  Price = "111.33";     Volume = "  -2.55%  " },
    @{ Row = 6;  Price = "263.99";     Volume = "  -1.77%  " },
    @{ Row = 7;  Price = $null;        Volume = "  +2.59%  " },
    @{ Row = 8;  Price = $null;        Volume = "  -0.24%  " },
    @{ Row = 9;  Price = "0.606";      Volume = "  -2.70%  " },
    @{ Row = 10; Price = "46.42";      Volume = "  -4.13%  " },
    @{ Row = 11; Price = "0.0935";     Volume = "  -1.83%  " },
    @{ Row = 12; Price = $null;        Volume = "  +1.96%  " },
    @{ Row = 13; Price = "0.109";      Volume = "  +1.74%  " },
    @{ Row = 14; Price = $null;        Volume = "  -2.50%  " },
    @{ Row = 15; Price = "2.614.37";   Volume = "  -0.98%  " },
    @{ Row = 16; Price = "0.862";      Volume = "  +1.29%  " },
    @{ Row = 17; Price = "2.263.82";   Volume = "  -1.44%  " },
    @{ Row = 18; Price = "43.092.01";  Volume = "  -1.49%  " },
    @{ Row = 19; Price = "0.0000108";  Volume = "  -2.36%  " },
    @{ Row = 20; Price = "6.73";       Volume = "  +1.47%  " },
    @{ Row = 21; Price = "71.89";      Volume = "  -0.95%  " },
    @{ Row = 22; Price = "2.44";       Volume = "  -1.97%  " },
    @{ Row = 23; Price = "233.54";     Volume = "  +0.10%  " },
    @{ Row = 26; Price = $null;        Volume = "  +1.99%  " },
    @{ Row = 27; Price = "11.28";      Volume = "  -3.47%  " },
    @{ Row = 28; Price = "41.22";      Volume = "  -1.71%  " },
    @{ Row = 29; Price = "3.34";       Volume = "  -1.71%  " },
    @{ Row = 30; Price = $null;        Volume = "  -0.41%  " },
    @{ Row = 31; Price = "173.23";     Volume = "  -2.09%  " },
    @{ Row = 32; Price = "21.39";      Volume = "  -1.39%  " },
    @{ Row = 33; Price = "0.0896";     Volume = "  -4.08%  " },
    @{ Row = 34; Price = "5.62";       Volume = "  -0.43%  " },
    @{ Row = 35; Price = $null;        Volume = "  +2.36%  " },
    @{ Row = 36; Price = "0.0373";     Volume = "  +2.52%  " },
    @{ Row = 37; Price = "4.63";       Volume = "  -2.86%  " },
    @{ Row = 38; Price = "3.92";       Volume = "  +2.55%  " },
    @{ Row = 39; Price = $null;        Volume = "  -3.99%  " },
    @{ Row = 40; Price = "2.58";       Volume = "  +7.15%  " },
    @{ Row = 41; Price = "14.30";      Volume = "  +2.29%  " },
    @{ Row = 42; Price = "75.72";      Volume = "  +5.11%  " },
    @{ Row = 43; Price = $null;        Volume = "  -4.22%  " },
    @{ Row = 44; Price = "6.08";       Volume = "  -2.09%  " },
    @{ Row = 45; Price = $null;        Volume = "  -0.04%  " },
    @{ Row = 46; Price = $null;        Volume = "  -4.08%  " },
    @{ Row = 47; Price = "8.51";       Volume = "  -3.63%  " },
    @{ Row = 50; Price = "100.51";     Volume = "  -2.01%  " }
)

foreach ($u in $priceVolumeUpdates) {
    if ($null -ne $u.Price) {
        Set-TextValue $ws.Cells.Item($u.Row, 4) $u.Price
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.Volume
}

# --- Rows whose contents (Coin / Link / Price / Volume) are fully replaced ---
# Row 24 and 25 swap ranking positions between PancakeSwap and
# InternetComputer(DFINITY); row 48/49 swap Cronos and TrustWalletToken;
# row 51 changes from TheSandbox to WOONetwork.
$fullRowUpdates = @(
    @{ Row = 24; Coin = "PancakeSwap";               Link = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake";                 Price = "2.85";     Volume = "  +0.72%  " },
    @{ Row = 25; Coin = "InternetComputer(DFINITY)";  Link = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp";       Price = "9.33";     Volume = "  -5.16%  " },
    @{ Row = 48; Coin = "TrustWalletToken";            Link = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";             Price = "1.25";     Volume = "  +1.43%  " },
    @{ Row = 49; Coin = "Cronos";                      Link = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";                    Price = "0.0994";   Volume = "  -1.52%  " },
    @{ Row = 51; Coin = "WOONetwork";                  Link = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo";                   Price = "0.432";    Volume = "  -3.91%  " }
)

foreach ($u in $fullRowUpdates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.Coin
    $ws.Cells.Item($u.Row, 3).Value = $u.Link
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.Price
    $ws.Cells.Item($u.Row, 5).Value = $u.Volume
}
